$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update values in column C (EARNED) for the two now-completed items
$ws.Range("C3").Value = 0.01
$ws.Range("C6").Value = 0.03

# Update the view: zoom normal to 100 (no custom zoomScale), selection moves to C4
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("C4").Select()
